$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.586.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.578"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.00%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "35.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.36%  "
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.081.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.840.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.648"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.589.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.59%  "
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "171.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.37%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.118"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0530"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.403.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.680"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.38%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.955"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0514"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.980.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0131"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.14%  "
